$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 297-298 (existing rows 297-350 shift down to 299-352),
# mirroring the weekly data-entry insert shown in the diff.
$ws.Rows("297:298").Insert()

# New row 297: Clementina / Especial, Region de O'Higgins, week of 45131
$ws.Range("A297").Value = 7
$ws.Range("B297").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C297").Value = "Ñuble"
$ws.Range("D297").Value = 45131
$ws.Range("E297").Value = 16
$ws.Range("F297").Value = "Fruta"
$ws.Range("G297").Value = 100102
$ws.Range("H297").Value = "Cítricos"
$ws.Range("I297").Value = 100102004
$ws.Range("J297").Value = "Mandarina"
$ws.Range("K297").Value = "Clementina"
$ws.Range("L297").Value = "Especial"
$ws.Range("M297").Value = 80
$ws.Range("N297").Value = 10000
$ws.Range("O297").Value = 10000
$ws.Range("P297").Value = 10000
$ws.Range("Q297").Value = "$/bandeja 10 kilos"
$ws.Range("R297").Value = "Región de O'Higgins"
$ws.Range("S297").Value = 1000
$ws.Range("T297").Value = 10

# New row 298: Clementina / Primera, Region de O'Higgins, week of 45131
$ws.Range("A298").Value = 7
$ws.Range("B298").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C298").Value = "Ñuble"
$ws.Range("D298").Value = 45131
$ws.Range("E298").Value = 16
$ws.Range("F298").Value = "Fruta"
$ws.Range("G298").Value = 100102
$ws.Range("H298").Value = "Cítricos"
$ws.Range("I298").Value = 100102004
$ws.Range("J298").Value = "Mandarina"
$ws.Range("K298").Value = "Clementina"
$ws.Range("L298").Value = "Primera"
$ws.Range("M298").Value = 120
$ws.Range("N298").Value = 9000
$ws.Range("O298").Value = 9000
$ws.Range("P298").Value = 9000
$ws.Range("Q298").Value = "$/bandeja 10 kilos"
$ws.Range("R298").Value = "Región de O'Higgins"
$ws.Range("S298").Value = 900
$ws.Range("T298").Value = 10
